$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.734583854675293
$ws.Range("B1").Value = 2.413678884506226
$ws.Range("C1").Value = 2.472666025161743
$ws.Range("D1").Value = 2.83130145072937
$ws.Range("E1").Value = 3.571537971496582
